# "Generate Report for Archive"
#
# 1) The localization status text changed from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F3, and the
#    "Status" column on the zh-cn / de-de detail sheets).
# 2) The "zh-cn"/"de-de" columns on the Overview sheet (E, F) and the
#    "Status" column (C) on the zh-cn / de-de detail sheets were narrowed,
#    presumably by an auto-fit pass re-run when the (now shorter) status
#    text changed.

$wb = $excel.ActiveWorkbook

# --- 1) Update the status text on every sheet -----------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2) Narrow the previously-widened columns ------------------------------
# ColumnWidth is rounded by the host to the nearest 1/6 character, so
# 12.5 is the value that lands on the narrow width the columns now use.
$newColumnWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E ("zh-cn")
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F ("de-de")

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C ("Status")

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C ("Status")
